# Template for "template_user.xlsx" was reworked: the two sample "Customer"
# rows were replaced by a single sample "member" row, and the header order
# was rearranged to put level_id first.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row (row 1): level_id, username, nama, password
$ws.Range("A1").Value = "level_id"
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "nama"
$ws.Range("D1").Value = "password"

# Sample data row (row 2): replaces the old "Customer 1" row with a new
# "member1" / "Khrisna" sample row, level_id changed from 5 to 6.
$ws.Range("A2").Value = 6
$ws.Range("B2").Value = "member1"
$ws.Range("C2").Value = "Khrisna"
$ws.Range("D2").Value = 12345

# The old row 3 ("Customer 2" / "customer2" sample) is removed entirely.
$ws.Range("A3:D3").ClearContents()
$ws.Rows("3:3").Delete()

# Restore the selection to A2 (was E8 before the edit).
$ws.Range("A2").Select() | Out-Null
